$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.728445529937744
$ws.Range("B1").Value = 3.94737982749939
$ws.Range("C1").Value = 3.478374242782593
$ws.Range("D1").Value = 4.362921714782715
$ws.Range("E1").Value = 5.182995319366455
